$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Heure " column (column C) entirely, shifting Killer and Points left
$ws.Range("C1:C2").EntireColumn.Delete()

# Update the selection to match the target state
$ws.Range("E2").Select()
